$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "62.875.83"
$ws.Cells.Item(2, 5).Value = "  -0.03%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.377.96"
$ws.Cells.Item(3, 5).Value = "  +0.58%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.01%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "562.51"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.53%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "154.50"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.71%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.00"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +0.16%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "3.374.35"
$ws.Cells.Item(8, 5).Value = "  +0.35%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.541"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +1.98%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  -2.17%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +1.88%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -1.69%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "3.960.81"
$ws.Cells.Item(13, 5).Value = "  +0.51%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  -3.71%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  +3.21%  "

# Row 16
$ws.Cells.Item(16, 5).Value = "  -1.12%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "62.956.64"
$ws.Cells.Item(17, 5).Value = "  +0.00%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "3.374.62"
$ws.Cells.Item(18, 5).Value = "  +1.81%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  -4.46%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  +0.94%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "376.55"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -3.45%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  -4.80%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -0.17%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "71.22"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.90%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.526"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -2.88%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.0000116"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +19.29%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "9.43"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +6.32%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  -0.29%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "5.98"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +5.25%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -0.87%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "RenderToken"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "6.42"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -4.01%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "Fetch.AI"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.33"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +1.37%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "EthereumClassic"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "22.94"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.46%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "USDe"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.998"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.01%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -0.21%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "157.55"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -1.96%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.0758"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +1.98%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "2.895.89"
$ws.Cells.Item(40, 5).Value = "  +2.35%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "EnergySwap"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "26.68"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -1.46%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "Stacks"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.80"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -4.75%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.0314"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +1.56%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "41.07"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.68%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "4.30"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.66%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.749"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -0.08%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "22.95"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +3.21%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +0.47%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  +15.83%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "SuiNetwork"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.829"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +2.66%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "Cosmos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "6.33"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.08%  "
